$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1726078799249531
$ws.Range("C2").Value = 0.6097560975609756
$ws.Range("J2").Value = 0.01125703564727955
$ws.Range("P2").Value = 0.1219512195121951
$ws.Range("S2").Value = 0.08442776735459662
$ws.Range("B3").Value = 0.002906976744186046
$ws.Range("C3").Value = 0.04069767441860465
$ws.Range("J3").Value = 0.02906976744186046
$ws.Range("P3").Value = 0.7267441860465116
$ws.Range("S3").Value = 0.2005813953488372
$ws.Range("J4").Value = 0.05681818181818182
$ws.Range("O4").Value = 0.01136363636363636
$ws.Range("P4").Value = 0.7159090909090909
$ws.Range("S4").Value = 0.2159090909090909
$ws.Range("J5").Value = 0.1
$ws.Range("P5").Value = 0.6
$ws.Range("S5").Value = 0.3
$ws.Range("B6").Value = 0.04210526315789474
$ws.Range("D6").Value = 0.02105263157894737
$ws.Range("E6").Value = 0.00631578947368421
$ws.Range("F6").Value = 0.08210526315789474
$ws.Range("J6").Value = 0.248421052631579
$ws.Range("O6").Value = 0.01052631578947368
$ws.Range("Q6").Value = 0.1747368421052632
$ws.Range("R6").Value = 0.06736842105263158
$ws.Range("S6").Value = 0.3473684210526316
$ws.Range("B7").Value = 0.08755760368663594
$ws.Range("D7").Value = 0.01152073732718894
$ws.Range("E7").Value = 0.002304147465437788
$ws.Range("F7").Value = 0.06912442396313365
$ws.Range("J7").Value = 0.1267281105990783
$ws.Range("O7").Value = 0.01152073732718894
$ws.Range("Q7").Value = 0.1728110599078341
$ws.Range("R7").Value = 0.08525345622119816
$ws.Range("S7").Value = 0.4331797235023042
$ws.Range("B8").Value = 0.09820485744456177
$ws.Range("D8").Value = 0.02006335797254488
$ws.Range("E8").Value = 0.002111932418162619
$ws.Range("F8").Value = 0.06124604012671594
$ws.Range("J8").Value = 0.1077085533262936
$ws.Range("O8").Value = 0.02006335797254488
$ws.Range("Q8").Value = 0.1763463569165787
$ws.Range("R8").Value = 0.08342133051742344
$ws.Range("S8").Value = 0.4308342133051742
$ws.Range("B9").Value = 0.0774487471526196
$ws.Range("D9").Value = 0.02050113895216401
$ws.Range("F9").Value = 0.07289293849658314
$ws.Range("J9").Value = 0.1070615034168565
$ws.Range("O9").Value = 0.01366742596810934
$ws.Range("Q9").Value = 0.1845102505694761
$ws.Range("R9").Value = 0.08200455580865604
$ws.Range("S9").Value = 0.4419134396355353
$ws.Range("B10").Value = 0.1022076860179885
$ws.Range("D10").Value = 0.01880621422730989
$ws.Range("E10").Value = 0.002044153720359771
$ws.Range("F10").Value = 0.07972199509403106
$ws.Range("J10").Value = 0.1116107931316435
$ws.Range("O10").Value = 0.01308258381030254
$ws.Range("Q10").Value = 0.2105478331970564
$ws.Range("R10").Value = 0.0874897792313982
$ws.Range("S10").Value = 0.3744889615699101
$ws.Range("G11").Value = 0.1453125
$ws.Range("J11").Value = 0.09218750000000001
$ws.Range("K11").Value = 0.20625
$ws.Range("L11").Value = 0.5453125
$ws.Range("S11").Value = 0.0109375
$ws.Range("G12").Value = 0.7837078651685393
$ws.Range("J12").Value = 0.1601123595505618
$ws.Range("K12").Value = 0.005617977528089887
$ws.Range("L12").Value = 0.01966292134831461
$ws.Range("S12").Value = 0.03089887640449438
$ws.Range("G13").Value = 0.6779661016949152
$ws.Range("J13").Value = 0.3050847457627119
$ws.Range("S13").Value = 0.01694915254237288
$ws.Range("F15").Value = 0.02272727272727273
$ws.Range("H15").Value = 0.2
$ws.Range("I15").Value = 0.07272727272727272
$ws.Range("J15").Value = 0.3431818181818182
$ws.Range("K15").Value = 0.06363636363636363
$ws.Range("M15").Value = 0.02045454545454545
$ws.Range("N15").Value = 0.002272727272727273
$ws.Range("O15").Value = 0.0659090909090909
$ws.Range("S15").Value = 0.2090909090909091
$ws.Range("F16").Value = 0.01333333333333333
$ws.Range("H16").Value = 0.2106666666666667
$ws.Range("I16").Value = 0.09866666666666667
$ws.Range("J16").Value = 0.3786666666666667
$ws.Range("K16").Value = 0.09066666666666667
$ws.Range("M16").Value = 0.024
$ws.Range("N16").Value = 0.002666666666666667
$ws.Range("O16").Value = 0.05333333333333334
$ws.Range("S16").Value = 0.128
$ws.Range("F17").Value = 0.02380952380952381
$ws.Range("H17").Value = 0.1764069264069264
$ws.Range("I17").Value = 0.09956709956709957
$ws.Range("J17").Value = 0.4025974025974026
$ws.Range("K17").Value = 0.1168831168831169
$ws.Range("M17").Value = 0.02705627705627706
$ws.Range("O17").Value = 0.06926406926406926
$ws.Range("S17").Value = 0.08441558441558442
$ws.Range("F18").Value = 0.01776649746192894
$ws.Range("H18").Value = 0.1776649746192893
$ws.Range("I18").Value = 0.116751269035533
$ws.Range("J18").Value = 0.383248730964467
$ws.Range("K18").Value = 0.1015228426395939
$ws.Range("M18").Value = 0.02791878172588833
$ws.Range("N18").Value = 0.002538071065989848
$ws.Range("O18").Value = 0.06598984771573604
$ws.Range("S18").Value = 0.1065989847715736
$ws.Range("F19").Value = 0.01367108966626458
$ws.Range("H19").Value = 0.2219541616405308
$ws.Range("I19").Value = 0.09529553679131483
$ws.Range("J19").Value = 0.3590671491757137
$ws.Range("K19").Value = 0.1178126256533977
$ws.Range("M19").Value = 0.02694008845999196
$ws.Range("N19").Value = 0.00120627261761158
$ws.Range("O19").Value = 0.06996381182147166
$ws.Range("S19").Value = 0.09408926417370325
